# Insert a new weekly data row at row 212 (pushing existing rows 212:338
# down to 213:339), then populate the new row with the latest week's
# Brócoli price record for "Macroferia Regional de Talca".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 212 downwards to make room for the new record.
$ws.Rows("212:212").Insert()

# Fill in the newly inserted row with the new week's data.
$ws.Range("A212").Value = 5
$ws.Range("B212").Value = "Macroferia Regional de Talca"
$ws.Range("C212").Value = "Maule"
$ws.Range("D212").Value = 44719
$ws.Range("E212").Value = 7
$ws.Range("F212").Value = 100112023
$ws.Range("G212").Value = "Brócoli"
$ws.Range("H212").Value = "Sin especificar"
$ws.Range("I212").Value = "Primera"
$ws.Range("J212").Value = 3000
$ws.Range("K212").Value = 1000
$ws.Range("L212").Value = 1000
$ws.Range("M212").Value = 1000
$ws.Range("N212").Value = "$/unidad"
$ws.Range("O212").Value = "Región del Maule"
$ws.Range("P212").Value = 1000
$ws.Range("Q212").Value = 1
$ws.Range("R212").Value = "Hortaliza"
